$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 and Row 17 swap their "Periodo Mora" (E) and "Valor Mora" (F) values.
# Before: E16=1907, F16=55120 ; E17=1906, F17=20211
# After:  E16=1906, F16=20211 ; E17=1907, F17=55120
$ws.Range("E16").Value = "1906"
$ws.Range("F16").Value = 20211
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 55120
